$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 21 de Agosto de 2020 a las 12:19"

# Row 4: Estados Unidos
$ws.Cells.Item(4,1).Value = "Estados Unidos"
$ws.Cells.Item(4,2).Value = 5746534
$ws.Cells.Item(4,3).Value = 262
$ws.Cells.Item(4,4).Value = 3095910
$ws.Cells.Item(4,5).Value = 2473186
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 14
$ws.Cells.Item(4,8).Value = 177438

# Row 14: Iran
$ws.Cells.Item(14,1).Value = "Iran"
$ws.Cells.Item(14,2).Value = 354764
$ws.Cells.Item(14,3).Value = 2206
$ws.Cells.Item(14,4).Value = 305866
$ws.Cells.Item(14,5).Value = 28522
$ws.Cells.Item(14,6).Value = 0
$ws.Cells.Item(14,7).Value = 112
$ws.Cells.Item(14,8).Value = 20376

# Row 19: Banglades
$ws.Cells.Item(19,1).Value = "Banglades"
$ws.Cells.Item(19,2).Value = 290360
$ws.Cells.Item(19,3).Value = 2401
$ws.Cells.Item(19,4).Value = 172615
$ws.Cells.Item(19,5).Value = 113884
$ws.Cells.Item(19,6).Value = 0
$ws.Cells.Item(19,7).Value = 39
$ws.Cells.Item(19,8).Value = 3861

# Row 25: Filipinas
$ws.Cells.Item(25,1).Value = "Filipinas"
$ws.Cells.Item(25,2).Value = 182365
$ws.Cells.Item(25,3).Value = 4786
$ws.Cells.Item(25,4).Value = 114519
$ws.Cells.Item(25,5).Value = 64906
$ws.Cells.Item(25,6).Value = 0
$ws.Cells.Item(25,7).Value = 59
$ws.Cells.Item(25,8).Value = 2940

# Row 26: Indonesia
$ws.Cells.Item(26,1).Value = "Indonesia"
$ws.Cells.Item(26,2).Value = 149408
$ws.Cells.Item(26,3).Value = 2197
$ws.Cells.Item(26,4).Value = 102991
$ws.Cells.Item(26,5).Value = 39917
$ws.Cells.Item(26,6).Value = 0
$ws.Cells.Item(26,7).Value = 82
$ws.Cells.Item(26,8).Value = 6500

# Row 41: Kuwait
$ws.Cells.Item(41,1).Value = "Kuwait"
$ws.Cells.Item(41,2).Value = 79269
$ws.Cells.Item(41,3).Value = 502
$ws.Cells.Item(41,4).Value = 71264
$ws.Cells.Item(41,5).Value = 7494
$ws.Cells.Item(41,6).Value = 0
$ws.Cells.Item(41,7).Value = 2
$ws.Cells.Item(41,8).Value = 511

# Row 47: Polonia
$ws.Cells.Item(47,1).Value = "Polonia"
$ws.Cells.Item(47,2).Value = 60281
$ws.Cells.Item(47,3).Value = 903
$ws.Cells.Item(47,4).Value = 41029
$ws.Cells.Item(47,5).Value = 17314
$ws.Cells.Item(47,6).Value = 0
$ws.Cells.Item(47,7).Value = 13
$ws.Cells.Item(47,8).Value = 1938

# Row 49: Singapur
$ws.Cells.Item(49,1).Value = "Singapur"
$ws.Cells.Item(49,2).Value = 56216
$ws.Cells.Item(49,3).Value = 117
$ws.Cells.Item(49,4).Value = 53119
$ws.Cells.Item(49,5).Value = 3070
$ws.Cells.Item(49,6).Value = 0
$ws.Cells.Item(49,7).Value = 0
$ws.Cells.Item(49,8).Value = 27

# Row 71: Austria
$ws.Cells.Item(71,1).Value = "Austria"
$ws.Cells.Item(71,2).Value = 24762
$ws.Cells.Item(71,3).Value = 331
$ws.Cells.Item(71,4).Value = 21260
$ws.Cells.Item(71,5).Value = 2772
$ws.Cells.Item(71,6).Value = 0
$ws.Cells.Item(71,7).Value = 1
$ws.Cells.Item(71,8).Value = 730

# Row 73: El Salvador
$ws.Cells.Item(73,1).Value = "El Salvador"
$ws.Cells.Item(73,2).Value = 24200
$ws.Cells.Item(73,3).Value = 236
$ws.Cells.Item(73,4).Value = 11781
$ws.Cells.Item(73,5).Value = 11773
$ws.Cells.Item(73,6).Value = 0
$ws.Cells.Item(73,7).Value = 6
$ws.Cells.Item(73,8).Value = 646

# Row 90: Consejo Danes para los Refugiados
$ws.Cells.Item(90,1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(90,2).Value = 9802
$ws.Cells.Item(90,3).Value = 45
$ws.Cells.Item(90,4).Value = 8920
$ws.Cells.Item(90,5).Value = 634
$ws.Cells.Item(90,6).Value = 0
$ws.Cells.Item(90,7).Value = 1
$ws.Cells.Item(90,8).Value = 248

# Row 92: Malasia
$ws.Cells.Item(92,1).Value = "Malasia"
$ws.Cells.Item(92,2).Value = 9249
$ws.Cells.Item(92,3).Value = 9
$ws.Cells.Item(92,4).Value = 8945
$ws.Cells.Item(92,5).Value = 179
$ws.Cells.Item(92,6).Value = 0
$ws.Cells.Item(92,7).Value = 0
$ws.Cells.Item(92,8).Value = 125

# Row 100: Finlandia
$ws.Cells.Item(100,1).Value = "Finlandia"
$ws.Cells.Item(100,2).Value = 7871
$ws.Cells.Item(100,3).Value = 29
$ws.Cells.Item(100,4).Value = 7100
$ws.Cells.Item(100,5).Value = 437
$ws.Cells.Item(100,6).Value = 0
$ws.Cells.Item(100,7).Value = 0
$ws.Cells.Item(100,8).Value = 334

# Row 104: Maldivas
$ws.Cells.Item(104,1).Value = "Maldivas"
$ws.Cells.Item(104,2).Value = 6370
$ws.Cells.Item(104,3).Value = 0
$ws.Cells.Item(104,4).Value = 3915
$ws.Cells.Item(104,5).Value = 2430
$ws.Cells.Item(104,6).Value = 0
$ws.Cells.Item(104,7).Value = 1
$ws.Cells.Item(104,8).Value = 25

# Row 112: Hong Kong
$ws.Cells.Item(112,1).Value = "Hong Kong"
$ws.Cells.Item(112,2).Value = 4632
$ws.Cells.Item(112,3).Value = 27
$ws.Cells.Item(112,4).Value = 3900
$ws.Cells.Item(112,5).Value = 657
$ws.Cells.Item(112,6).Value = 0
$ws.Cells.Item(112,7).Value = 2
$ws.Cells.Item(112,8).Value = 75

# Row 122: Eslovaquia
$ws.Cells.Item(122,1).Value = "Eslovaquia"
$ws.Cells.Item(122,2).Value = 3225
$ws.Cells.Item(122,3).Value = 123
$ws.Cells.Item(122,4).Value = 2045
$ws.Cells.Item(122,5).Value = 1147
$ws.Cells.Item(122,6).Value = 0
$ws.Cells.Item(122,7).Value = 0
$ws.Cells.Item(122,8).Value = 33

# Row 123: Mayotte
$ws.Cells.Item(123,1).Value = "Mayotte"
$ws.Cells.Item(123,2).Value = 3160
$ws.Cells.Item(123,3).Value = 0
$ws.Cells.Item(123,4).Value = 2964
$ws.Cells.Item(123,5).Value = 157
$ws.Cells.Item(123,6).Value = 0
$ws.Cells.Item(123,7).Value = 0
$ws.Cells.Item(123,8).Value = 39

# Row 124: Mozambique
$ws.Cells.Item(124,1).Value = "Mozambique"
$ws.Cells.Item(124,2).Value = 3115
$ws.Cells.Item(124,3).Value = 0
$ws.Cells.Item(124,4).Value = 1380
$ws.Cells.Item(124,5).Value = 1715
$ws.Cells.Item(124,6).Value = 0
$ws.Cells.Item(124,7).Value = 0
$ws.Cells.Item(124,8).Value = 20

# Row 125: Sri Lanka
$ws.Cells.Item(125,1).Value = "Sri Lanka"
$ws.Cells.Item(125,2).Value = 2918
$ws.Cells.Item(125,3).Value = 0
$ws.Cells.Item(125,4).Value = 2789
$ws.Cells.Item(125,5).Value = 118
$ws.Cells.Item(125,6).Value = 0
$ws.Cells.Item(125,7).Value = 0
$ws.Cells.Item(125,8).Value = 11

# Row 128: Eslovenia
$ws.Cells.Item(128,1).Value = "Eslovenia"
$ws.Cells.Item(128,2).Value = 2574
$ws.Cells.Item(128,3).Value = 38
$ws.Cells.Item(128,4).Value = 2079
$ws.Cells.Item(128,5).Value = 365
$ws.Cells.Item(128,6).Value = 0
$ws.Cells.Item(128,7).Value = 1
$ws.Cells.Item(128,8).Value = 130

# Row 129: Lituania
$ws.Cells.Item(129,1).Value = "Lituania"
$ws.Cells.Item(129,2).Value = 2564
$ws.Cells.Item(129,3).Value = 36
$ws.Cells.Item(129,4).Value = 1755
$ws.Cells.Item(129,5).Value = 726
$ws.Cells.Item(129,6).Value = 0
$ws.Cells.Item(129,7).Value = 1
$ws.Cells.Item(129,8).Value = 83

# Row 130: Tunez
$ws.Cells.Item(130,1).Value = "Tunez"
$ws.Cells.Item(130,2).Value = 2543
$ws.Cells.Item(130,3).Value = 0
$ws.Cells.Item(130,4).Value = 1397
$ws.Cells.Item(130,5).Value = 1083
$ws.Cells.Item(130,6).Value = 0
$ws.Cells.Item(130,7).Value = 0
$ws.Cells.Item(130,8).Value = 63

Write-Host "Update complete"
